$d = $word.ActiveDocument

# --- Update the letter date: "Le 29 octobre 2025" -> "Le 5 novembre 2025" ---
$d.Content.Find.Execute("Le 29 octobre 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Le 5 novembre 2025", 2)

# --- Re-touch the "Objet" line so the two runs collapse back into one
#     (matches how Word re-serializes the paragraph after an in-place edit) ---
$d.Content.Find.Execute(" : Postulation pour une demande de stage en développement application web",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         " : Postulation pour une demande de stage en développement application web", 2)

# --- Re-touch the skills lines so the "<label> :" / " <value>" runs collapse
#     back into a single run (and the stale gramEnd proofing marker between
#     them is dropped), leaving the bold "Frontend"/"Backend" label run intact ---
$d.Content.Find.Execute(" : JavaScript, HTML, CSS, Tailwind CSS",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         " : JavaScript, HTML, CSS, Tailwind CSS", 2)

$d.Content.Find.Execute(" : PHP, MySQL, PostgreSQL, Node.js, Express.js, WebSocket",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         " : PHP, MySQL, PostgreSQL, Node.js, Express.js, WebSocket", 2)
